$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27 (shifts existing rows 27-70 down to 28-71)
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new record's data
$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value = "Ñuble"
$ws.Cells.Item(27, 4).Value = 44571
$ws.Cells.Item(27, 5).Value = 16
$ws.Cells.Item(27, 6).Value = 100112030
$ws.Cells.Item(27, 7).Value = "Poroto granado"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 120
$ws.Cells.Item(27, 11).Value = 29000
$ws.Cells.Item(27, 12).Value = 30000
$ws.Cells.Item(27, 13).Value = 29500
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(27, 16).Value = 1180
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
